$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for Price column (D) so numeric-looking strings
# like "1.033" or "27.735.25" are stored as text, matching the source data.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.735.25"
$ws.Range("D3").Value = "1.864.39"
$ws.Range("D5").Value = "324.22"
$ws.Range("D6").Value = "1.033"
$ws.Range("D7").Value = "0.4413"
$ws.Range("D8").Value = "0.3800"
$ws.Range("D9").Value = "0.07456"
$ws.Range("D10").Value = "0.8845"
$ws.Range("D12").Value = "1.867.81"
$ws.Range("D13").Value = "5.560"
$ws.Range("D14").Value = "6.748"
$ws.Range("D15").Value = "0.07200"
$ws.Range("D16").Value = "84.01"
$ws.Range("D18").Value = "0.000009104"
$ws.Range("D19").Value = "1.033"
$ws.Range("D20").Value = "15.54"
$ws.Range("D21").Value = "27.756.87"
$ws.Range("D22").Value = "5.312"
$ws.Range("D24").Value = "158.52"
$ws.Range("D25").Value = "1.945"
$ws.Range("D26").Value = "18.83"
$ws.Range("D27").Value = "1.991"
$ws.Range("D28").Value = "5.314"
$ws.Range("D29").Value = "117.74"
$ws.Range("D30").Value = "0.09100"
$ws.Range("D31").Value = "1.215"
$ws.Range("D32").Value = "0.7715"
$ws.Range("D33").Value = "3.010"
$ws.Range("D34").Value = "4.575"
$ws.Range("D35").Value = "1.034"
$ws.Range("D36").Value = "1.163"
$ws.Range("D38").Value = "0.05348"
$ws.Range("D39").Value = "2.847"
$ws.Range("D40").Value = "0.5197"
$ws.Range("D41").Value = "0.1694"
$ws.Range("D42").Value = "6.856"
$ws.Range("D43").Value = "8.715"
$ws.Range("D44").Value = "109.58"
$ws.Range("D45").Value = "10.61"
$ws.Range("D46").Value = "1.732"
$ws.Range("D47").Value = "0.4697"
$ws.Range("D48").Value = "0.06428"
$ws.Range("D49").Value = "1.875"
$ws.Range("D50").Value = "39.74"
$ws.Range("D51").Value = "64.50"

# Restore default (General) style now that text values are committed.
$ws.Range("D2:D51").Style = "Normal"

# Coin name / link / volume columns are plain text and percent strings;
# these do not get mis-parsed as numbers, so assign directly.
$ws.Range("E2").Value = "  +3.07%  "
$ws.Range("E3").Value = "  +3.01%  "
$ws.Range("E4").Value = "  +2.81%  "
$ws.Range("E5").Value = "  +3.95%  "
$ws.Range("E6").Value = "  +2.63%  "
$ws.Range("E7").Value = "  +2.89%  "
$ws.Range("E8").Value = "  +3.10%  "
$ws.Range("E9").Value = "  +2.75%  "
$ws.Range("E10").Value = "  +2.30%  "
$ws.Range("E11").Value = "  +3.01%  "
$ws.Range("E12").Value = "  -8.82%  "
$ws.Range("E13").Value = "  +3.08%  "
$ws.Range("E14").Value = "  +1.92%  "
$ws.Range("E15").Value = "  +3.90%  "
$ws.Range("E16").Value = "  +4.04%  "
$ws.Range("E17").Value = "  +2.57%  "
$ws.Range("E18").Value = "  +2.93%  "
$ws.Range("E19").Value = "  +2.64%  "
$ws.Range("E20").Value = "  +2.42%  "
$ws.Range("E21").Value = "  +3.00%  "
$ws.Range("E22").Value = "  +2.29%  "
$ws.Range("E23").Value = "  +4.63%  "
$ws.Range("E24").Value = "  +2.93%  "
$ws.Range("E25").Value = "  +3.31%  "
$ws.Range("E27").Value = "  +4.08%  "
$ws.Range("E28").Value = "  +1.50%  "
$ws.Range("E29").Value = "  +2.67%  "
$ws.Range("E30").Value = "  +1.73%  "
$ws.Range("E31").Value = "  +4.96%  "
$ws.Range("E32").Value = "  +4.13%  "
$ws.Range("E33").Value = "  +7.30%  "
$ws.Range("E34").Value = "  +3.20%  "
$ws.Range("E35").Value = "  +2.72%  "
$ws.Range("E36").Value = "  +4.39%  "
$ws.Range("E37").Value = "  +3.31%  "
$ws.Range("E38").Value = "  +2.49%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("E39").Value = "  +3.16%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("E40").Value = "  +2.27%  "
$ws.Range("E41").Value = "  +2.81%  "
$ws.Range("E42").Value = "  +6.13%  "
$ws.Range("E43").Value = "  +5.27%  "
$ws.Range("E44").Value = "  +2.13%  "
$ws.Range("E45").Value = "  +2.22%  "
$ws.Range("E46").Value = "  +5.33%  "
$ws.Range("E47").Value = "  +2.73%  "
$ws.Range("E48").Value = "  +2.49%  "
$ws.Range("E49").Value = "  +3.38%  "
$ws.Range("E50").Value = "  +4.72%  "
$ws.Range("E51").Value = "  +1.58%  "
